# Add a new "Supervisors" sheet, positioned before the existing "Projects"
# sheet (Worksheets.Add() without args inserts at the front, matching the
# target tab order: Supervisors, Projects, Student_preferences,
# Supervisor_preferences).
$wb = $excel.ActiveWorkbook

$supervisors = $wb.Worksheets.Add()
$supervisors.Name = "Supervisors"

$supervisors.Range("A1").Value = "Supervisor"
$supervisors.Range("B1").Value = "Max_number_of_projects"
$supervisors.Range("C1").Value = "Max_number_of_students"
$supervisors.Range("A2").Value = "Dr Smith"

# Update the "Projects" sheet: rename the first column header, and add a
# "Supervisor" column naming the (single, so far) supervisor for every
# existing project.
$projects = $wb.Worksheets.Item("Projects")
$projects.Range("A1").Value = "Project"
$projects.Range("C1").Value = "Supervisor"
$projects.Range("C2").Value = "Dr Smith"
$projects.Range("C3").Value = "Dr Smith"
$projects.Range("C4").Value = "Dr Smith"
$projects.Range("C5").Value = "Dr Smith"
$projects.Range("C6").Value = "Dr Smith"

$wb.Worksheets.Item(1).Activate()
